$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6542786666666667
$ws.Range("H2").Value = 1.962836
$ws.Range("I2").Value = 0.3193006097963691
$ws.Range("J2").Value = 0.3193006097963691
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 5.771161376390222
$ws.Range("R2").Value = 51.940452387512
$ws.Range("S2").Value = 0.02048360240602649
$ws.Range("T2").Value = 0.02048360240602649
# Row 3
$ws.Range("G3").Value = 0.6542786666666667
$ws.Range("H3").Value = 1.962836
$ws.Range("I3").Value = 0.3193006097963691
$ws.Range("J3").Value = 0.3193006097963691
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 35.79659017108666
$ws.Range("R3").Value = 322.16931153978
$ws.Range("S3").Value = 0.1270529574091806
$ws.Range("T3").Value = 0.1270529574091806
# Row 4
$ws.Range("G4").Value = 0.6542786666666667
$ws.Range("H4").Value = 1.962836
$ws.Range("I4").Value = 0.3193006097963691
$ws.Range("J4").Value = 0.3193006097963691
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 14.33404389484889
$ws.Range("R4").Value = 129.00639505364
$ws.Range("S4").Value = 0.05087587001357891
$ws.Range("T4").Value = 0.05087587001357891
# Row 5
$ws.Range("G5").Value = 0.6542786666666667
$ws.Range("H5").Value = 1.962836
$ws.Range("I5").Value = 0.3193006097963691
$ws.Range("J5").Value = 0.3193006097963691
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 34.05969229737467
$ws.Range("R5").Value = 306.537230676372
$ws.Range("S5").Value = 0.1208881799675831
$ws.Range("T5").Value = 0.1208881799675831
# Row 6
$ws.Range("I6").Value = 0.4124821994964292
$ws.Range("J6").Value = 0.4124821994964292
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 7.455361077138
$ws.Range("R6").Value = 67.098249694242
$ws.Range("S6").Value = 0.02646133804578858
$ws.Range("T6").Value = 0.02646133804578858
# Row 7
$ws.Range("I7").Value = 0.4124821994964292
$ws.Range("J7").Value = 0.4124821994964292
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.1641308588733578
$ws.Range("T7").Value = 0.1641308588733578
# Row 8
$ws.Range("I8").Value = 0.4124821994964292
$ws.Range("J8").Value = 0.4124821994964292
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 18.51715208811
$ws.Range("R8").Value = 166.65436879299
$ws.Range("S8").Value = 0.06572298993690831
$ws.Range("T8").Value = 0.06572298993690831
# Row 9
$ws.Range("I9").Value = 0.4124821994964292
$ws.Range("J9").Value = 0.4124821994964292
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 43.999342193403
$ws.Range("R9").Value = 395.9940797406269
$ws.Range("S9").Value = 0.1561670126403744
$ws.Range("T9").Value = 0.1561670126403744
# Row 10
$ws.Range("G10").Value = 0.5380133333333333
$ws.Range("H10").Value = 1.61404
$ws.Range("I10").Value = 0.2625608844731457
$ws.Range("J10").Value = 0.2625608844731457
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 4.745625873964444
$ws.Range("R10").Value = 42.71063286568
$ws.Range("S10").Value = 0.01684366581182686
$ws.Range("T10").Value = 0.01684366581182686
# Row 11
$ws.Range("G11").Value = 0.5380133333333333
$ws.Range("H11").Value = 1.61404
$ws.Range("I11").Value = 0.2625608844731457
$ws.Range("J11").Value = 0.2625608844731457
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 29.43553531713333
$ws.Range("R11").Value = 264.9198178541999
$ws.Range("S11").Value = 0.1044756441071561
$ws.Range("T11").Value = 0.1044756441071561
# Row 12
$ws.Range("G12").Value = 0.5380133333333333
$ws.Range("H12").Value = 1.61404
$ws.Range("I12").Value = 0.2625608844731457
$ws.Range("J12").Value = 0.2625608844731457
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 11.78688398217778
$ws.Range("R12").Value = 106.0819558396
$ws.Range("S12").Value = 0.04183522680280822
$ws.Range("T12").Value = 0.04183522680280822
# Row 13
$ws.Range("G13").Value = 0.5380133333333333
$ws.Range("H13").Value = 1.61404
$ws.Range("I13").Value = 0.2625608844731457
$ws.Range("J13").Value = 0.2625608844731457
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 28.00728423345333
$ws.Range("R13").Value = 252.0655581010799
$ws.Range("S13").Value = 0.09940634775135454
$ws.Range("T13").Value = 0.09940634775135455
# Row 14
$ws.Range("G14").Value = 0.01159033333333333
$ws.Range("H14").Value = 0.034771
$ws.Range("I14").Value = 0.005656306234056004
$ws.Range("J14").Value = 0.005656306234056004
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 0.1022342428091111
$ws.Range("R14").Value = 0.9201081852820001
$ws.Range("S14").Value = 0.0003628603404767118
$ws.Range("T14").Value = 0.0003628603404767118
# Row 15
$ws.Range("G15").Value = 0.01159033333333333
$ws.Range("H15").Value = 0.034771
$ws.Range("I15").Value = 0.005656306234056004
$ws.Range("J15").Value = 0.005656306234056004
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 0.6341249278283334
$ws.Range("R15").Value = 5.707124350455
$ws.Range("S15").Value = 0.00225070173059523
$ws.Range("T15").Value = 0.00225070173059523
# Row 16
$ws.Range("G16").Value = 0.01159033333333333
$ws.Range("H16").Value = 0.034771
$ws.Range("I16").Value = 0.005656306234056004
$ws.Range("J16").Value = 0.005656306234056004
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 0.2539229157544445
$ws.Range("R16").Value = 2.28530624179
$ws.Range("S16").Value = 0.0009012494555032371
$ws.Range("T16").Value = 0.0009012494555032371
# Row 17
$ws.Range("G17").Value = 0.01159033333333333
$ws.Range("H17").Value = 0.034771
$ws.Range("I17").Value = 0.005656306234056004
$ws.Range("J17").Value = 0.005656306234056004
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 0.6033563480963333
$ws.Range("R17").Value = 5.430207132866999
$ws.Range("S17").Value = 0.002141494707480824
$ws.Range("T17").Value = 0.002141494707480824
